# Update countries & provincias Spain
#
# 1) Two pairs of countries were reordered in the source data (their rank
#    rows keep the same statistics columns, but the country name shown on
#    that rank changes), and 2) a batch of rows received refreshed
#    case/death counts for the current data pull.
#
# Row 34 (rank 38) now shows "Rumania" (used to be "Corea del Sur")
# Row 35 (rank 39) now shows "Corea del Sur" (used to be "Rumania")
# Row 36 (rank 40) now shows "Bielorrusia" (used to be "Emiratos Arabes Unidos")
# Row 37 (rank 41) now shows "Emiratos Arabes Unidos" (used to be "Bielorrusia")
# Row 151 (rank 155) now shows "Sierra Leona" (used to be "Zambia")
# Row 152 (rank 156) now shows "Zambia" (used to be "Sierra Leona")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: Rumania (new data)
$ws.Range("A34").Value = "Rumania"
$ws.Range("B34").Value = 11036
$ws.Range("C34").Value = 401
$ws.Range("D34").Value = 3054
$ws.Range("E34").Value = 7374
$ws.Range("F34").Value = 236
$ws.Range("G34").Value = 7
$ws.Range("H34").Value = 608

# Row 35: Corea del Sur (carries the figures that used to sit on row 34)
$ws.Range("A35").Value = "Corea del Sur"
$ws.Range("B35").Value = 10728
$ws.Range("C35").Value = 10
$ws.Range("D35").Value = 8717
$ws.Range("E35").Value = 1769
$ws.Range("F35").Value = 55
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 242

# Row 36: Bielorrusia (new data)
$ws.Range("A36").Value = "Bielorrusia"
$ws.Range("B36").Value = 10463
$ws.Range("C36").Value = 873
$ws.Range("D36").Value = 1695
$ws.Range("E36").Value = 8696
$ws.Range("F36").Value = 92
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 72

# Row 37: Emiratos Arabes Unidos (carries the figures that used to sit on row 36)
$ws.Range("A37").Value = "Emiratos Arabes Unidos"
$ws.Range("B37").Value = 9813
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 1887
$ws.Range("E37").Value = 7855
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 71

# Row 55: Marruecos - refreshed figures
$ws.Range("B55").Value = 4047
$ws.Range("C55").Value = 150
$ws.Range("D55").Value = 557
$ws.Range("E55").Value = 3330
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 160

# Row 75: Bosnia y Herzegovina - refreshed figures
$ws.Range("B75").Value = 1516
$ws.Range("C75").Value = 30
$ws.Range("D75").Value = 624
$ws.Range("E75").Value = 833
$ws.Range("F75").Value = 4
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 59

# Row 97: Libano - refreshed figures
$ws.Range("B97").Value = 707
$ws.Range("C97").Value = 3
$ws.Range("D97").Value = 145
$ws.Range("E97").Value = 538
$ws.Range("F97").Value = 44
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 24

# Row 109: Malta - refreshed figures
$ws.Range("B109").Value = 448
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 282
$ws.Range("E109").Value = 162
$ws.Range("F109").Value = 2
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 4

# Row 151: Sierra Leona (new data)
$ws.Range("A151").Value = "Sierra Leona"
$ws.Range("B151").Value = 86
$ws.Range("C151").Value = 4
$ws.Range("D151").Value = 10
$ws.Range("E151").Value = 73
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = 3

# Row 152: Zambia (carries the figures that used to sit on row 151)
$ws.Range("A152").Value = "Zambia"
$ws.Range("B152").Value = 84
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 37
$ws.Range("E152").Value = 44
$ws.Range("F152").Value = 1
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 3

# Row 169: Macao - refreshed figures
$ws.Range("B169").Value = 45
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 31
$ws.Range("E169").Value = 14
$ws.Range("F169").Value = 1
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 0
